$d = $word.ActiveDocument
$enDash = [char]0x2013

# --- 1) Merge "...Siddhi " + "(" runs into a single run (text unchanged). ---
$d.Content.Find.Execute(
    "Drug Name, Drug Inchikey, Drug Targets (semi-colon separated) $enDash Siddhi (",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Drug Name, Drug Inchikey, Drug Targets (semi-colon separated) $enDash Siddhi (",
    2) | Out-Null

# --- 2) Merge "...in a sample – " + "Raghvendra - " runs into a single run. ---
$d.Content.Find.Execute(
    "f) Convert the mutation information into a matrix where each sample is row and columns are list of genes and each entry is number of times a gene was mutated in a sample $enDash Raghvendra - ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "f) Convert the mutation information into a matrix where each sample is row and columns are list of genes and each entry is number of times a gene was mutated in a sample $enDash Raghvendra - ",
    2) | Out-Null

# --- 3) Merge "...for each sample- " + "Raghvendra - " runs into a single run. ---
$d.Content.Find.Execute(
    "g) Get the list of marker genes for celltype scores and estimate celltype score for each sample- Raghvendra - ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "g) Get the list of marker genes for celltype scores and estimate celltype score for each sample- Raghvendra - ",
    2) | Out-Null

# --- 4) Change the hyphen before "Raghvendra" to an en dash on the "I)" line. ---
$d.Content.Find.Execute(
    "I) Get the pathway AUC score for each cell line, pathway combined - Raghvendra",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I) Get the pathway AUC score for each cell line, pathway combined $enDash Raghvendra",
    2) | Out-Null

# --- 5) Append new paragraphs (two blank ones, a heading, and items a-d). ---
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()

$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()

$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$new = $d.Paragraphs($d.Paragraphs.Count)
$new.Range.Text = "Action Points (21.05.2023):"

$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$new = $d.Paragraphs($d.Paragraphs.Count)
$new.Range.Text = "a) Generate the mapping with <sample id, module class>. - From Raghvendra to Siddhi"

$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$new = $d.Paragraphs($d.Paragraphs.Count)
$new.Range.Text = "b) Identify the optimal set of varying genes using k-means/t-sne with the module class to quantitatively identify the set $enDash From Siddhi to Raghvendra"

$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$new = $d.Paragraphs($d.Paragraphs.Count)
$new.Range.Text = "c) Generate the training and test set for ML models by taking union of oncogenes and varying genes $enDash From Raghvendra"

$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$new = $d.Paragraphs($d.Paragraphs.Count)
$new.Range.Text = "d) Divide the task of building ML models with this new train/test set - Raghvendra/Siddhi"
